# This script rewrites the worker account-statement rows (B16:J50) of Hoja1.
# The previous account statement rows are replaced by a new set of rows for
# 5 workers x 7 periods (2502..2408), matching the refreshed database export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# JAVIER JESUS VARGAS PEREZ
$ws.Range("C16").Value = "1047390509"
$ws.Range("D16").Value = "JAVIER JESUS VARGAS PEREZ"
$ws.Range("E16").Value = "2502"
$ws.Range("F16").Value = 32933
$ws.Range("C17").Value = "1047390509"
$ws.Range("D17").Value = "JAVIER JESUS VARGAS PEREZ"
$ws.Range("E17").Value = "2501"
$ws.Range("F17").Value = 52000
$ws.Range("C18").Value = "1047390509"
$ws.Range("D18").Value = "JAVIER JESUS VARGAS PEREZ"
$ws.Range("E18").Value = "2412"
$ws.Range("F18").Value = 52000
$ws.Range("C19").Value = "1047390509"
$ws.Range("D19").Value = "JAVIER JESUS VARGAS PEREZ"
$ws.Range("E19").Value = "2411"
$ws.Range("F19").Value = 52000
$ws.Range("C20").Value = "1047390509"
$ws.Range("D20").Value = "JAVIER JESUS VARGAS PEREZ"
$ws.Range("E20").Value = "2410"
$ws.Range("F20").Value = 52000
$ws.Range("C21").Value = "1047390509"
$ws.Range("D21").Value = "JAVIER JESUS VARGAS PEREZ"
$ws.Range("E21").Value = "2409"
$ws.Range("F21").Value = 52000
$ws.Range("C22").Value = "1047390509"
$ws.Range("D22").Value = "JAVIER JESUS VARGAS PEREZ"
$ws.Range("E22").Value = "2408"
$ws.Range("F22").Value = 45066

# ALONSO JOSE VELASCO CARRILLO
$ws.Range("C23").Value = "1127584034"
$ws.Range("D23").Value = "ALONSO JOSE VELASCO CARRILLO"
$ws.Range("E23").Value = "2502"
$ws.Range("F23").Value = 32933
$ws.Range("C24").Value = "1127584034"
$ws.Range("D24").Value = "ALONSO JOSE VELASCO CARRILLO"
$ws.Range("E24").Value = "2501"
$ws.Range("F24").Value = 52000
$ws.Range("C25").Value = "1127584034"
$ws.Range("D25").Value = "ALONSO JOSE VELASCO CARRILLO"
$ws.Range("E25").Value = "2412"
$ws.Range("F25").Value = 52000
$ws.Range("C26").Value = "1127584034"
$ws.Range("D26").Value = "ALONSO JOSE VELASCO CARRILLO"
$ws.Range("E26").Value = "2411"
$ws.Range("F26").Value = 52000
$ws.Range("C27").Value = "1127584034"
$ws.Range("D27").Value = "ALONSO JOSE VELASCO CARRILLO"
$ws.Range("E27").Value = "2410"
$ws.Range("F27").Value = 52000
$ws.Range("C28").Value = "1127584034"
$ws.Range("D28").Value = "ALONSO JOSE VELASCO CARRILLO"
$ws.Range("E28").Value = "2409"
$ws.Range("F28").Value = 52000
$ws.Range("C29").Value = "1127584034"
$ws.Range("D29").Value = "ALONSO JOSE VELASCO CARRILLO"
$ws.Range("E29").Value = "2408"
$ws.Range("F29").Value = 52000

# JUAN CAMILO LOPEZ RIOS
$ws.Range("C30").Value = "1103117470"
$ws.Range("D30").Value = "JUAN CAMILO LOPEZ RIOS"
$ws.Range("E30").Value = "2502"
$ws.Range("F30").Value = 32933
$ws.Range("C31").Value = "1103117470"
$ws.Range("D31").Value = "JUAN CAMILO LOPEZ RIOS"
$ws.Range("E31").Value = "2501"
$ws.Range("F31").Value = 52000
$ws.Range("C32").Value = "1103117470"
$ws.Range("D32").Value = "JUAN CAMILO LOPEZ RIOS"
$ws.Range("E32").Value = "2412"
$ws.Range("F32").Value = 52000
$ws.Range("C33").Value = "1103117470"
$ws.Range("D33").Value = "JUAN CAMILO LOPEZ RIOS"
$ws.Range("E33").Value = "2411"
$ws.Range("F33").Value = 52000
$ws.Range("C34").Value = "1103117470"
$ws.Range("D34").Value = "JUAN CAMILO LOPEZ RIOS"
$ws.Range("E34").Value = "2410"
$ws.Range("F34").Value = 52000
$ws.Range("C35").Value = "1103117470"
$ws.Range("D35").Value = "JUAN CAMILO LOPEZ RIOS"
$ws.Range("E35").Value = "2409"
$ws.Range("F35").Value = 52000
$ws.Range("C36").Value = "1103117470"
$ws.Range("D36").Value = "JUAN CAMILO LOPEZ RIOS"
$ws.Range("E36").Value = "2408"
$ws.Range("F36").Value = 52000

# RAUL RUIZ RAMOS
$ws.Range("C37").Value = "1085038750"
$ws.Range("D37").Value = "RAUL RUIZ RAMOS"
$ws.Range("E37").Value = "2502"
$ws.Range("F37").Value = 32933
$ws.Range("C38").Value = "1085038750"
$ws.Range("D38").Value = "RAUL RUIZ RAMOS"
$ws.Range("E38").Value = "2501"
$ws.Range("F38").Value = 52000
$ws.Range("C39").Value = "1085038750"
$ws.Range("D39").Value = "RAUL RUIZ RAMOS"
$ws.Range("E39").Value = "2412"
$ws.Range("F39").Value = 52000
$ws.Range("C40").Value = "1085038750"
$ws.Range("D40").Value = "RAUL RUIZ RAMOS"
$ws.Range("E40").Value = "2411"
$ws.Range("F40").Value = 52000
$ws.Range("C41").Value = "1085038750"
$ws.Range("D41").Value = "RAUL RUIZ RAMOS"
$ws.Range("E41").Value = "2410"
$ws.Range("F41").Value = 52000
$ws.Range("C42").Value = "1085038750"
$ws.Range("D42").Value = "RAUL RUIZ RAMOS"
$ws.Range("E42").Value = "2409"
$ws.Range("F42").Value = 52000
$ws.Range("C43").Value = "1085038750"
$ws.Range("D43").Value = "RAUL RUIZ RAMOS"
$ws.Range("E43").Value = "2408"
$ws.Range("F43").Value = 52000

# OLGA ISABEL VASQUEZ TEHERAN
$ws.Range("C44").Value = "1052071317"
$ws.Range("D44").Value = "OLGA ISABEL VASQUEZ TEHERAN"
$ws.Range("E44").Value = "2502"
$ws.Range("F44").Value = 32933
$ws.Range("C45").Value = "1052071317"
$ws.Range("D45").Value = "OLGA ISABEL VASQUEZ TEHERAN"
$ws.Range("E45").Value = "2501"
$ws.Range("F45").Value = 52000
$ws.Range("C46").Value = "1052071317"
$ws.Range("D46").Value = "OLGA ISABEL VASQUEZ TEHERAN"
$ws.Range("E46").Value = "2412"
$ws.Range("F46").Value = 52000
$ws.Range("C47").Value = "1052071317"
$ws.Range("D47").Value = "OLGA ISABEL VASQUEZ TEHERAN"
$ws.Range("E47").Value = "2411"
$ws.Range("F47").Value = 52000
$ws.Range("C48").Value = "1052071317"
$ws.Range("D48").Value = "OLGA ISABEL VASQUEZ TEHERAN"
$ws.Range("E48").Value = "2410"
$ws.Range("F48").Value = 52000
$ws.Range("C49").Value = "1052071317"
$ws.Range("D49").Value = "OLGA ISABEL VASQUEZ TEHERAN"
$ws.Range("E49").Value = "2409"
$ws.Range("F49").Value = 52000
$ws.Range("C50").Value = "1052071317"
$ws.Range("D50").Value = "OLGA ISABEL VASQUEZ TEHERAN"
$ws.Range("E50").Value = "2408"
$ws.Range("F50").Value = 45066
